# Portfolio header rename: "DEPOSIT" -> "DEPOSITS"
# (and the derived product_description labels that embed it)
# on the PORTFOLIO_SNAP sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PORTFOLIO_SNAP")

# A single whole-word-insensitive substring replace over the used range
# turns:
#   "DEPOSIT"                   -> "DEPOSITS"                   (product_category header, column I)
#   "Fixed Deposit"             -> "Fixed DEPOSITS"              (product_description, column E)
#   "Recurring Deposit Account" -> "Recurring DEPOSITS Account"  (product_description, column E)
# in one pass, matching every affected cell.
$ws.Cells.Replace("DEPOSIT", "DEPOSITS")
